# The original document has two paragraphs:
#   1. "Food words"
#   2. an otherwise-empty paragraph that only carries the <w:bookmarkStart>/
#      <w:bookmarkEnd> pair for the "_GoBack" bookmark (Word silently drops
#      this bookmark at whatever point the cursor last sat).
#
# This edit fleshes the page out with the body copy for two recipe write-ups.
# The text was pasted in from a blog/CMS and kept its raw markup as literal
# characters -- <h2>, <p>, <br>, <img> -- rather than being turned into Word
# formatting, exactly as the diff shows (everything is plain text, so those
# angle brackets show up XML-escaped as &lt;/&gt; once saved). A few trailing
# blank paragraphs are left at the end -- the author's commit message notes
# the inline pictures meant to go there, with text wrapping, never got
# finished ("tried to text wrap around pictures but not working just yet").
#
# The "_GoBack" bookmark has to stay exactly where it is, so instead of
# clobbering that paragraph we keep its bookmarkStart/bookmarkEnd pair in
# place and grow all of the new paragraphs around it with a single
# Range.InsertXML call. The bookmark ends up sitting mid-paragraph in
# "It started with boiling ... bubbled a<bookmark/>way while we cooked ...",
# i.e. right where the cursor was when the author kept typing through it.

$d = $word.ActiveDocument

# The bookmark-only paragraph is the last paragraph in the (as yet almost
# empty) document.
$bookmarkPara = $d.Paragraphs($d.Paragraphs.Count)
$r = $bookmarkPara.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>&lt;h2&gt;Sausage and Fennel Pasta&lt;/h2&gt;</w:t></w:r></w:p><w:p><w:r><w:t>&lt;p&gt;This one was a winter warmer from Jamie Oliver. Simple to make, yummy to eat, and according to the dashing Essex Chef, a ‘superfood’ meal – what’s not to love? &lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>br</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;</w:t></w:r></w:p><w:p><w:r><w:t>It started with boiling some chipolata sausages, broccoli stalk and chili pepper for a few minutes. Then everything is fried off with onions, garlic, and fennel seeds and oregano. A can of tinned tomatoes was then added and everything bubbled a</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>way while we cooked the pasta. The broccoli florets were added to cook with the pasta in the last few minutes then everything was brought together in the pan. &lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>br</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">We ‘paired’ this dish with a wine from Tesco that was on special. It was a Californian red and the label recommended to </w:t></w:r><w:r><w:t>drink</w:t></w:r><w:r><w:t xml:space="preserve"> it with meaty pasta, so it had to be. </w:t></w:r><w:r><w:t>&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>br</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;</w:t></w:r></w:p><w:p><w:r><w:t>Just as we were finishing dinner, a cry of ‘I can make apple sauce!’ came from Nathalie. True to her word, twenty minutes later a delicious plate of apple sauce topped with a gingery-oat crumble, raspberry, and (not for Sam) a dollop of plain yoghurt arrived. The perfect light dessert for a hearty meal.</w:t></w:r></w:p><w:p><w:r><w:t>&lt;/p&gt;</w:t></w:r></w:p><w:p><w:r><w:t>&lt;h2&gt;Chicken Pesto and Veg&lt;/h2&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">&lt;p&gt;Impromptu dinner parties are a favourite at 7D, and this meal was one of them. </w:t></w:r><w:r><w:t>A last-</w:t></w:r><w:r><w:t>minute throw together meal, inspired from what was in the 7D fridge and an Aldi run</w:t></w:r><w:r><w:t xml:space="preserve"> courtesy of an AU rental van. </w:t></w:r><w:r><w:t>&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>br</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;</w:t></w:r></w:p><w:p><w:r><w:t>The conversation for this to happen happened as follows:</w:t></w:r></w:p><w:p><w:r><w:t>&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>img</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt; &lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>img</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;</w:t></w:r></w:p><w:p><w:r><w:t>Sam already had chicken and broccoli in the fridge, so Rebecca picked up some pesto, and ingredients for a light apple, pear, and avocado salad. The dressing was a mustard and balsamic number, which Laura even ventured to have! Major win. &lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>br</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;</w:t></w:r></w:p><w:p><w:r><w:t>&lt;/p&gt;</w:t></w:r></w:p><w:p><w:r><w:t>&lt;h2&gt;</w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p/>'

# Splice the whole run of new paragraphs into that range. InsertXML parses
# the WordprocessingML and inserts it in place of the (essentially empty)
# range, preserving the bookmark pair that is embedded part-way through.
$r.InsertXML($xml)
